$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 136.23077
$ws.Range("I33").Value = 148.27272
$ws.Range("J33").Value = 70
$ws.Range("K33").Value = 148.27272
$ws.Range("L33").Value = 70
$ws.Range("M33").Value = 80.72728000000001
$ws.Range("N33").Value = -528

$ws.Range("H38").Value = 690.8889
$ws.Range("I38").Value = 36.333332
$ws.Range("J38").Value = 2000
$ws.Range("K38").Value = 108.999996
$ws.Range("L38").Value = 6000
$ws.Range("M38").Value = 263.000004
$ws.Range("N38").Value = -6744

$ws.Range("H40").Value = 6533.7354
$ws.Range("I40").Value = 1735
$ws.Range("J40").Value = 8533.208000000001
$ws.Range("K40").Value = 1735
$ws.Range("L40").Value = 8533.208000000001
$ws.Range("M40").Value = -1560
$ws.Range("N40").Value = -8883.208000000001

$ws.Range("H70").Value = 1019.8
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 1024.75
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 3074.25
$ws.Range("M70").Value = -2730
$ws.Range("N70").Value = -3614.25

$ws.Range("H73").Value = 1019.8
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 1024.75
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 3074.25
$ws.Range("M73").Value = -2064
$ws.Range("N73").Value = -4946.25

$ws.Range("H132").Value = 2028
$ws.Range("I132").Value = 929
$ws.Range("J132").Value = 4226
$ws.Range("K132").Value = 2787
$ws.Range("L132").Value = 12678
$ws.Range("M132").Value = -257
$ws.Range("N132").Value = -17738

$ws.Range("H138").Value = 3322.9285
$ws.Range("I138").Value = 844
$ws.Range("J138").Value = 3999
$ws.Range("K138").Value = 2532
$ws.Range("L138").Value = 11997
$ws.Range("M138").Value = 2608
$ws.Range("N138").Value = -22277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 9996.272000000001
$ws.Range("I122").Value = 9996.299999999999
$ws.Range("J122").Value = 9996
$ws.Range("K122").Value = 29988.9
$ws.Range("L122").Value = 29988
$ws.Range("M122").Value = -27538.9
$ws.Range("N122").Value = -34888

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1642.6
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 1553.25
$ws.Range("K68").Value = 6000
$ws.Range("L68").Value = 4659.75
$ws.Range("M68").Value = -5189
$ws.Range("N68").Value = -6281.75

$ws.Range("H69").Value = 212
$ws.Range("I69").Value = 212
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 636
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 175

$ws.Range("H71").Value = 1642.6
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 1553.25
$ws.Range("K71").Value = 18000
$ws.Range("L71").Value = 13979.25
$ws.Range("M71").Value = -13944
$ws.Range("N71").Value = -22091.25

$ws.Range("H72").Value = 212
$ws.Range("I72").Value = 212
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 1908
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = 2148

$ws.Range("H86").Value = 1518.75
$ws.Range("I86").Value = 358.33334
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 1075.00002
$ws.Range("L86").Value = 15000
$ws.Range("M86").Value = 110.9999800000001
$ws.Range("N86").Value = -17372

$ws.Range("H89").Value = 1518.75
$ws.Range("I89").Value = 358.33334
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 3225.00006
$ws.Range("L89").Value = 45000
$ws.Range("M89").Value = 2702.99994
$ws.Range("N89").Value = -56856

$ws.Range("H97").Value = 316
$ws.Range("I97").Value = 316
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 948
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -452
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()

$ws.Range("H7").Value = 25002500
$ws.Range("I7").Value = 25002500
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 25002500
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -25002388
$ws.Range("N7").ClearContents()

$ws.Range("H8").Value = 25002500
$ws.Range("I8").Value = 25002500
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 25002500
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -25002361
$ws.Range("N8").ClearContents()

$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()

$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

$ws.Range("H21").Value = 5666666.5
$ws.Range("I21").Value = 5666666.5
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 5666666.5
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -5666493.5

$ws.Range("H25").Value = 200
$ws.Range("I25").Value = 200
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 200
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 329

$ws.Range("H30").Value = 5666666.5
$ws.Range("I30").Value = 5666666.5
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 5666666.5
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -5666561.5

$ws.Range("H31").Value = 5155
$ws.Range("I31").Value = 5155
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 5155
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -4863

$ws.Range("H36").Value = 5001320
$ws.Range("I36").Value = 5001320
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 5001320
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -5000835
$ws.Range("N36").ClearContents()

$ws.Range("H37").Value = 5155
$ws.Range("I37").Value = 5155
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 5155
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -4878

$ws.Range("H43").Value = 7291.2856
$ws.Range("I43").Value = 1509.75
$ws.Range("J43").Value = 15000
$ws.Range("K43").Value = 1509.75
$ws.Range("L43").Value = 15000
$ws.Range("M43").Value = -1358.75
$ws.Range("N43").Value = -15302

$ws.Range("H49").Value = 30000
$ws.Range("I49").Value = 10000
$ws.Range("J49").Value = 50000
$ws.Range("K49").Value = 10000
$ws.Range("L49").Value = 50000
$ws.Range("M49").Value = -9816
$ws.Range("N49").Value = -50368

$ws.Range("H55").Value = 19000
$ws.Range("I55").Value = 19000
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 19000
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -18673
$ws.Range("N55").ClearContents()

$ws.Range("H107").Value = 364.9091
$ws.Range("I107").Value = 301.1
$ws.Range("J107").Value = 1003
$ws.Range("K107").Value = 301.1
$ws.Range("L107").Value = 1003
$ws.Range("M107").Value = 1618.9
$ws.Range("N107").Value = -4843

$ws.Range("H122").Value = 10600
$ws.Range("I122").Value = 10600
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 31800
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -29350

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 405
$ws.Range("I18").Value = 405
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 405
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -233

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()

$ws.Range("H45").Value = 5000
$ws.Range("I45").Value = 5000
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 5000
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -4593
$ws.Range("N45").ClearContents()

$ws.Range("H46").Value = 934
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 1002
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 1002
$ws.Range("M46").Value = -712
$ws.Range("N46").Value = -1378

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8500
$ws.Range("I122").Value = 8500
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 25500
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -23050
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 3833.3333
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -18440
